# Refresh the cryptocurrency table (columns: B=Coin, C=Link, D=Price,
# E=Volume(1h)) with the latest snapshot values, per commit message:
# "Updated cryptos list on Tue Jul  9 16:13:36 UTC 2024 with GitHub Actions"
#
# Column D ("Price") holds plain-text figures in the source sheet, including
# ones that look numeric ("25.60", "1.00", "0.373", ...). Assigning those as a
# bare .Value lets Excel auto-convert them to real numbers, which would silently
# drop significant trailing zeros (25.60 -> 25.6, 1.00 -> 1). To keep them as
# text - matching the source data - every Price write is prefixed with a leading
# apostrophe, exactly like typing "'25.60" into a cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'57.536.27"
$ws.Cells.Item(2, 5).Value = "  +2.75%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'3.068.87"
$ws.Cells.Item(3, 5).Value = "  +2.78%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.01%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'516.86"
$ws.Cells.Item(5, 5).Value = "  +2.74%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'140.48"
$ws.Cells.Item(6, 5).Value = "  +1.71%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.45%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'7.27"
$ws.Cells.Item(9, 5).Value = "  +2.06%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.70%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.373"
$ws.Cells.Item(11, 5).Value = "  +2.76%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'3.594.29"
$ws.Cells.Item(12, 5).Value = "  +2.81%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +3.27%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'25.60"
$ws.Cells.Item(14, 5).Value = "  -1.55%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +2.07%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'57.589.09"
$ws.Cells.Item(16, 5).Value = "  +2.72%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'6.15"
$ws.Cells.Item(17, 5).Value = "  +2.33%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'3.062.22"
$ws.Cells.Item(18, 5).Value = "  +2.72%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'12.91"
$ws.Cells.Item(19, 5).Value = "  -0.27%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'8.11"
$ws.Cells.Item(20, 5).Value = "  +1.65%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'333.50"
$ws.Cells.Item(21, 5).Value = "  +1.75%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.03%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +1.40%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'65.93"
$ws.Cells.Item(24, 5).Value = "  +2.05%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +5.17%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.00%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'0.0₃0908"
$ws.Cells.Item(27, 5).Value = "  +1.16%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'6.38"
$ws.Cells.Item(28, 5).Value = "  -0.20%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'7.19"
$ws.Cells.Item(29, 5).Value = "  +3.08%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +2.89%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'20.77"
$ws.Cells.Item(31, 5).Value = "  +2.80%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'1.19"
$ws.Cells.Item(32, 5).Value = "  +2.96%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'155.13"
$ws.Cells.Item(33, 5).Value = "  +0.96%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'27.35"
$ws.Cells.Item(34, 5).Value = "  +8.20%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.94%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'5.94"
$ws.Cells.Item(36, 5).Value = "  +3.89%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'1.26"
$ws.Cells.Item(37, 5).Value = "  +2.02%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.0674"
$ws.Cells.Item(38, 5).Value = "  +2.56%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'3.108.12"
$ws.Cells.Item(39, 5).Value = "  +2.83%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'3.91"
$ws.Cells.Item(40, 5).Value = "  +3.63%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'36.92"
$ws.Cells.Item(41, 5).Value = "  +0.82%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'1.00"
$ws.Cells.Item(42, 5).Value = "  -0.03%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.656"
$ws.Cells.Item(43, 5).Value = "  +1.17%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "VeChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(44, 4).Value = "'0.0259"
$ws.Cells.Item(44, 5).Value = "  +9.69%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Maker"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(45, 4).Value = "'2.257.81"
$ws.Cells.Item(45, 5).Value = "  +4.55%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'1.37"
$ws.Cells.Item(46, 5).Value = "  +1.92%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'20.60"
$ws.Cells.Item(47, 5).Value = "  +4.84%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +1.06%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.34%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'266.30"
$ws.Cells.Item(50, 5).Value = "  +17.25%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "SuiNetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(51, 4).Value = "'0.706"
$ws.Cells.Item(51, 5).Value = "  +5.14%  "

Write-Output "Applied cryptos list update."
